$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# ---------------------------------------------------------------
# 1) Rename the "fedcore" column header to "approach" everywhere.
# ---------------------------------------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------
# 2) Give the empty header cells that sit inside the merged
#    B1:D1 / E1:G1 ranges a border that frames the whole merged
#    block instead of boxing every individual cell:
#      - the "inner" cell (C1 / F1) gets a top+bottom border only
#      - the "outer right" cell (D1 / G1) gets a top+bottom+right
#        border (closing the box on the right-hand side)
# ---------------------------------------------------------------
$r1 = $ws1.Range("C1:D1")
$r1.ClearFormats()
$r1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$r1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$r1.Borders.Item(10).LineStyle = 1  # xlEdgeRight (only reaches the D1 edge)

$r2 = $ws2.Range("C1:D1")
$r2.ClearFormats()
$r2.Borders.Item(8).LineStyle = 1
$r2.Borders.Item(9).LineStyle = 1
$r2.Borders.Item(10).LineStyle = 1

$r3 = $ws2.Range("F1:G1")
$r3.ClearFormats()
$r3.Borders.Item(8).LineStyle = 1
$r3.Borders.Item(9).LineStyle = 1
$r3.Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------
# 3) Drop the stray empty cell G5 on computational_comparison -
#    it used to hold an empty inline string and should no longer
#    be present at all.
# ---------------------------------------------------------------
$ws2.Range("G5").ClearContents()

Write-Host "edits applied"
